$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark after "ramp rates may be too fast."
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {}

# 2. Simplify the "Bridging pins 32 and 33..." paragraph - remove proofErr
#    markers and merge into a single run by replacing the full sentence text.
$d.Content.Find.Execute("Bridging pins 32 and 33, and changing MODE logic (software) from RB4 (pin 33) to RA8 (pin 32) corrects the problem.", $true, $false, $false, $false, $false, $true, 1, $false, "Bridging pins 32 and 33, and changing MODE logic (software) from RB4 (pin 33) to RA8 (pin 32) corrects the problem.", 2) | Out-Null

# 3. Append the new PCB bring-up test-log section (page break, header fields,
#    results table) at the end of the document, followed by the _GoBack
#    bookmark paragraph.
$full = $d.Content
$endRange = $d.Range($full.End, $full.End)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:br w:type="page"/></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>HW part number: HW00</w:t></w:r><w:r><w:t>01, HW0002</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">HW revision: </w:t></w:r><w:r><w:t xml:space="preserve">HW0001 </w:t></w:r><w:r><w:t>R002</w:t></w:r><w:r><w:t>, HW0002 R001</w:t></w:r></w:p><w:p><w:r><w:t>Firmware version: initial, unreleased test firmware</w:t></w:r></w:p><w:p><w:r><w:t>Software version: N/A</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Test case: board </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bringup</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>Test date: 2018-08-31</w:t></w:r></w:p><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="8095"/><w:gridCol w:w="1255"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="8095" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Test power supplies</w:t></w:r></w:p><w:p><w:r><w:t>Plug in USB mini connector to apply +5V.</w:t></w:r></w:p><w:p><w:r><w:t>Verify 5.0V appears at 5.0V power rail.</w:t></w:r></w:p><w:p><w:r><w:t>Verify 3.3V appears at 3.3V power rail (HW0001 R002)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1255" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Pass</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Pass</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="8095" w:type="dxa"/></w:tcPr><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Investigative notes</w:t></w:r></w:p><w:p><w:r><w:t>Power up sequence</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">0 – 3.3 V in 990 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>usec</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>Waveshape irregular with several dips as voltage ramps on. Voltage ramp during startup exceeds specification, but appears to work.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">PIC24FJ128GB204 datasheet calls out 0-3.3V in 66 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ms.</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:r><w:t>Recommended action: redesign logic board to control the ramp rate.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1255" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Design</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr></w:p></w:tc></w:tr></w:tbl><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endRange.InsertXML($xml)

Write-Output "done"
